$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    "55+7=",
    "99-91=",
    "81-28=",
    "43+15=",
    "31+52=",
    "40+39=",
    "16+57=",
    "85-59=",
    "90-76=",
    "96-37=",
    "77+6=",
    "4+63=",
    "81-5=",
    "0+52=",
    "44-3=",
    "0+83=",
    "14+45=",
    "31+48=",
    "28+21=",
    "83-78=",
    "11+63=",
    "10-1=",
    "83+3=",
    "56+23=",
    "80-31=",
    "95-57=",
    "32+45=",
    "90-40=",
    "70-66=",
    "25+18=",
    "29+49=",
    "2+14=",
    "52+12=",
    "83-39=",
    "27+44=",
    "59-5=",
    "87-12=",
    "98-40=",
    "63-0=",
    "12+44=",
    "99-18=",
    "0+75=",
    "53+15=",
    "70-64=",
    "56-8=",
    "75-25=",
    "74-39=",
    "3+59=",
    "98-77=",
    "81-56=",
    "32-12=",
    "33-21=",
    "24+57=",
    "72-66=",
    "93-73=",
    "57+3=",
    "54+15=",
    "67-47=",
    "80-69=",
    "20-18=",
    "81+12=",
    "64+4=",
    "60+7=",
    "42-20=",
    "18+38=",
    "8-6=",
    "79+8=",
    "4+39=",
    "29+48=",
    "99-71=",
    "97-84=",
    "56+29=",
    "25+31=",
    "88+8=",
    "52+12=",
    "11+14=",
    "62-24=",
    "5-1=",
    "71-61=",
    "18+58=",
    "93-52=",
    "45+40=",
    "27+66=",
    "98-74=",
    "70-63=",
    "32+9=",
    "50-13=",
    "23-10=",
    "0+16=",
    "89+2=",
    "12+21=",
    "59-15=",
    "32-3=",
    "21-17=",
    "41+31=",
    "95-58=",
    "63-18=",
    "12+80=",
    "13+39=",
    "52-32="
)

$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Replaced $idx cells"
